$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before existing row 191, pushing current rows 191-195 down to 194-198
$ws.Range("A191:A193").EntireRow.Insert()

# New row 191 - Valencia / Primera / Region de O'Higgins
$ws.Cells.Item(191, 1).Value2 = 1
$ws.Cells.Item(191, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(191, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(191, 4).Value2 = 45239
$ws.Cells.Item(191, 5).Value2 = 15
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value2 = 100102
$ws.Cells.Item(191, 8).Value = "Cítricos"
$ws.Cells.Item(191, 9).Value2 = 100102005
$ws.Cells.Item(191, 10).Value = "Naranja"
$ws.Cells.Item(191, 11).Value = "Valencia"
$ws.Cells.Item(191, 12).Value = "Primera"
$ws.Cells.Item(191, 13).Value2 = 200
$ws.Cells.Item(191, 14).Value2 = 850
$ws.Cells.Item(191, 15).Value2 = 900
$ws.Cells.Item(191, 16).Value2 = 875
$ws.Cells.Item(191, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(191, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(191, 19).Value2 = 875
$ws.Cells.Item(191, 20).Value2 = 1

# New row 192 - Valencia / Segunda / Region de O'Higgins
$ws.Cells.Item(192, 1).Value2 = 1
$ws.Cells.Item(192, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(192, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(192, 4).Value2 = 45239
$ws.Cells.Item(192, 5).Value2 = 15
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value2 = 100102
$ws.Cells.Item(192, 8).Value = "Cítricos"
$ws.Cells.Item(192, 9).Value2 = 100102005
$ws.Cells.Item(192, 10).Value = "Naranja"
$ws.Cells.Item(192, 11).Value = "Valencia"
$ws.Cells.Item(192, 12).Value = "Segunda"
$ws.Cells.Item(192, 13).Value2 = 270
$ws.Cells.Item(192, 14).Value2 = 750
$ws.Cells.Item(192, 15).Value2 = 800
$ws.Cells.Item(192, 16).Value2 = 775
$ws.Cells.Item(192, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(192, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(192, 19).Value2 = 775
$ws.Cells.Item(192, 20).Value2 = 1

# New row 193 - Valencia / Tercera / Region de O'Higgins
$ws.Cells.Item(193, 1).Value2 = 1
$ws.Cells.Item(193, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(193, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(193, 4).Value2 = 45239
$ws.Cells.Item(193, 5).Value2 = 15
$ws.Cells.Item(193, 6).Value = "Fruta"
$ws.Cells.Item(193, 7).Value2 = 100102
$ws.Cells.Item(193, 8).Value = "Cítricos"
$ws.Cells.Item(193, 9).Value2 = 100102005
$ws.Cells.Item(193, 10).Value = "Naranja"
$ws.Cells.Item(193, 11).Value = "Valencia"
$ws.Cells.Item(193, 12).Value = "Tercera"
$ws.Cells.Item(193, 13).Value2 = 300
$ws.Cells.Item(193, 14).Value2 = 650
$ws.Cells.Item(193, 15).Value2 = 700
$ws.Cells.Item(193, 16).Value2 = 675
$ws.Cells.Item(193, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(193, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(193, 19).Value2 = 675
$ws.Cells.Item(193, 20).Value2 = 1

"done"
